# Auto-generated script to apply cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '29.942.90'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '1.633.43'
$ws.Range("E3").Value = '  +1.76%  '

$ws.Range("E4").Value = '  -0.04%  '

Set-TextValue "D5" '214.43'
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("E7").Value = '  -0.09%  '

Set-TextValue "D8" '28.59'
$ws.Range("E8").Value = '  +1.75%  '

$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("E12").Value = '  +1.77%  '

$ws.Range("D13").Value = '1.629.85'
$ws.Range("E13").Value = '  +1.63%  '

Set-TextValue "D14" '0.563'
$ws.Range("E14").Value = '  +2.64%  '

Set-TextValue "D15" '9.27'
$ws.Range("E15").Value = '  +17.80%  '

$ws.Range("E16").Value = '  +2.62%  '

$ws.Range("D17").Value = '29.968.85'
$ws.Range("E17").Value = '  +0.91%  '

Set-TextValue "D18" '64.11'
$ws.Range("E18").Value = '  +0.01%  '

Set-TextValue "D19" '242.91'
$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("D20").Value = '0.0₃0701'
$ws.Range("E20").Value = '  +0.53%  '

Set-TextValue "D22" '9.86'
$ws.Range("E22").Value = '  +4.80%  '

$ws.Range("E23").Value = '  +2.42%  '

Set-TextValue "D24" '2.13'
$ws.Range("E24").Value = '  +1.15%  '

Set-TextValue "D25" '157.61'
$ws.Range("E25").Value = '  +1.42%  '

Set-TextValue "D26" '15.53'
$ws.Range("E26").Value = '  +0.41%  '

$ws.Range("E27").Value = '  +1.45%  '

Set-TextValue "D28" '6.61'
$ws.Range("E28").Value = '  +2.34%  '

$ws.Range("E29").Value = '  -0.05%  '

Set-TextValue "D30" '0.0485'
$ws.Range("E30").Value = '  +1.01%  '

$ws.Range("E31").Value = '  +4.09%  '

$ws.Range("E32").Value = '  +4.20%  '

$ws.Range("E33").Value = '  -0.63%  '

$ws.Range("D34").Value = '1.422.95'
$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("E35").Value = '  +4.61%  '

Set-TextValue "D36" '1.04'
$ws.Range("E36").Value = '  +0.41%  '

Set-TextValue "D37" '2.80'
$ws.Range("E37").Value = '  -3.85%  '

$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("E39").Value = '  +0.34%  '

Set-TextValue "D40" '75.94'
$ws.Range("E40").Value = '  +14.54%  '

Set-TextValue "D41" '0.553'
$ws.Range("E41").Value = '  +0.84%  '

Set-TextValue "D42" '1.99'
$ws.Range("E42").Value = '  +2.51%  '

Set-TextValue "D43" '0.829'
$ws.Range("E43").Value = '  +1.43%  '

Set-TextValue "D44" '0.0488'
$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D45" '0.999'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D46" '52.92'
$ws.Range("E46").Value = '  -6.85%  '

$ws.Range("E47").Value = '  +3.24%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D48" '5.35'
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '1.776.23'
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("E50").Value = '  +7.68%  '

Set-TextValue "D51" '89.55'
$ws.Range("E51").Value = '  +3.38%  '

